# Auto-generated Excel COM-interop script
# Applies numeric corrections to Sheets/Seraph_Profits.xlsx per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1000.5283
$ws.Range("I15").Value = 1000.5283
$ws.Range("K15").Value = 3001.5849
$ws.Range("M15").Value = -2832.5849
$ws.Range("H64").Value = 3206.25
$ws.Range("H67").Value = 3206.25
$ws.Range("H88").Value = 6750.5
$ws.Range("J88").Value = 6750.5
$ws.Range("L88").Value = 6750.5
$ws.Range("N88").Value = -7562.5
$ws.Range("H91").Value = 6750.5
$ws.Range("J91").Value = 6750.5
$ws.Range("L91").Value = 6750.5
$ws.Range("N91").Value = -9558.5
$ws.Range("H92").Value = 1141.8
$ws.Range("J92").Value = 1050
$ws.Range("L92").Value = 1050
$ws.Range("N92").Value = -3546
$ws.Range("H112").Value = 2043.6
$ws.Range("J112").Value = 2043.6
$ws.Range("L112").Value = 6130.799999999999
$ws.Range("N112").Value = -8346.799999999999
$ws.Range("H132").Value = 1804
$ws.Range("I132").Value = 1967
$ws.Range("K132").Value = 5901
$ws.Range("M132").Value = -3371
$ws.Range("H137").Value = 2741.889
$ws.Range("J137").Value = 3575
$ws.Range("L137").Value = 10725
$ws.Range("N137").Value = -15825

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2007.9412
$ws.Range("I74").Value = 1199.7858
$ws.Range("J74").Value = 5779.3335
$ws.Range("K74").Value = 1199.7858
$ws.Range("L74").Value = 5779.3335
$ws.Range("M74").Value = -325.7858000000001
$ws.Range("N74").Value = -7527.3335
$ws.Range("H77").Value = 2007.9412
$ws.Range("I77").Value = 1199.7858
$ws.Range("J77").Value = 5779.3335
$ws.Range("K77").Value = 5998.929
$ws.Range("L77").Value = 28896.6675
$ws.Range("M77").Value = -1630.929
$ws.Range("N77").Value = -37632.6675
$ws.Range("H132").Value = 1732.5714
$ws.Range("I132").Value = 1591.4546
$ws.Range("K132").Value = 4774.3638
$ws.Range("M132").Value = -2244.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 3800.6667
$ws.Range("I12").Value = 2900
$ws.Range("J12").Value = 4251
$ws.Range("K12").Value = 2900
$ws.Range("L12").Value = 4251
$ws.Range("M12").Value = -2732
$ws.Range("N12").Value = -4587
$ws.Range("H99").Value = 2007.2778
$ws.Range("I99").Value = 1991.1072
$ws.Range("J99").Value = 2063.875
$ws.Range("K99").Value = 1991.1072
$ws.Range("L99").Value = 2063.875
$ws.Range("M99").Value = -493.1071999999999
$ws.Range("N99").Value = -5059.875
$ws.Range("H105").Value = 2494.7273
$ws.Range("I105").Value = 2382.95
$ws.Range("J105").Value = 3612.5
$ws.Range("K105").Value = 2382.95
$ws.Range("L105").Value = 3612.5
$ws.Range("M105").Value = -635.9499999999998
$ws.Range("N105").Value = -7106.5
$ws.Range("H140").Value = 111030
$ws.Range("J140").Value = 111030
$ws.Range("L140").Value = 111030
$ws.Range("N140").Value = -121390

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 229.75
$ws.Range("I7").Value = 122.07692
$ws.Range("K7").Value = 122.07692
$ws.Range("M7").Value = -9.076920000000001
$ws.Range("H19").Value = 233.33333
$ws.Range("I19").Value = 233.33333
$ws.Range("K19").Value = 233.33333
$ws.Range("M19").Value = -63.33332999999999
$ws.Range("H24").Value = 233.33333
$ws.Range("I24").Value = 233.33333
$ws.Range("K24").Value = 233.33333
$ws.Range("M24").Value = -63.33332999999999
$ws.Range("H31").Value = 4163
$ws.Range("I31").Value = 3922.8667
$ws.Range("J31").Value = 4319.609
$ws.Range("K31").Value = 3922.8667
$ws.Range("L31").Value = 4319.609
$ws.Range("M31").Value = -3627.8667
$ws.Range("N31").Value = -4909.609
$ws.Range("H34").Value = 4163
$ws.Range("I34").Value = 3922.8667
$ws.Range("J34").Value = 4319.609
$ws.Range("K34").Value = 3922.8667
$ws.Range("L34").Value = 4319.609
$ws.Range("M34").Value = -3720.8667
$ws.Range("N34").Value = -4723.609
$ws.Range("H58").Value = 2788.1765
$ws.Range("I58").Value = 1265.45
$ws.Range("J58").Value = 4963.5
$ws.Range("K58").Value = 1265.45
$ws.Range("L58").Value = 4963.5
$ws.Range("M58").Value = -1062.45
$ws.Range("N58").Value = -5369.5
$ws.Range("H62").Value = 84017.8
$ws.Range("I62").Value = 5295
$ws.Range("J62").Value = 136499.67
$ws.Range("K62").Value = 5295
$ws.Range("L62").Value = 136499.67
$ws.Range("M62").Value = -4671
$ws.Range("N62").Value = -137747.67
$ws.Range("H65").Value = 84017.8
$ws.Range("I65").Value = 5295
$ws.Range("J65").Value = 136499.67
$ws.Range("K65").Value = 26475
$ws.Range("L65").Value = 682498.3500000001
$ws.Range("M65").Value = -23355
$ws.Range("N65").Value = -688738.3500000001
$ws.Range("H86").Value = 9857
$ws.Range("J86").Value = 11419.8
$ws.Range("L86").Value = 11419.8
$ws.Range("N86").Value = -13665.8
$ws.Range("H89").Value = 9857
$ws.Range("J89").Value = 11419.8
$ws.Range("L89").Value = 57099
$ws.Range("N89").Value = -68331
$ws.Range("H132").Value = 2064.4
$ws.Range("I132").Value = 1714.3334
$ws.Range("K132").Value = 5143.0002
$ws.Range("M132").Value = -2613.0002
$ws.Range("H136").Value = 2788.1765
$ws.Range("I136").Value = 1265.45
$ws.Range("J136").Value = 4963.5
$ws.Range("K136").Value = 3796.35
$ws.Range("L136").Value = 14890.5
$ws.Range("M136").Value = -1246.35
$ws.Range("N136").Value = -19990.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 120.375
$ws.Range("J38").Value = 192.25
$ws.Range("L38").Value = 576.75
$ws.Range("N38").Value = -1270.75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H121").Value = 1124.75
$ws.Range("I121").Value = 700
$ws.Range("J121").Value = 1266.3334
$ws.Range("K121").Value = 2100
$ws.Range("L121").Value = 3799.0002
$ws.Range("M121").Value = -790
$ws.Range("N121").Value = -6419.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 4999.3335
$ws.Range("J13").Value = 4999.3335
$ws.Range("L13").Value = 4999.3335
$ws.Range("N13").Value = -5277.3335
$ws.Range("H23").Value = 3500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3500
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -3946
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 55000
$ws.Range("L62").Value = 55000
$ws.Range("N62").Value = -56372
$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 55000
$ws.Range("L65").Value = 165000
$ws.Range("N65").Value = -171864
$ws.Range("H70").Value = 7394.3335
$ws.Range("I70").Value = 7312.25
$ws.Range("K70").Value = 7312.25
$ws.Range("M70").Value = -7042.25
$ws.Range("H73").Value = 7394.3335
$ws.Range("I73").Value = 7312.25
$ws.Range("K73").Value = 7312.25
$ws.Range("M73").Value = -6376.25
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 4753.25
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 2200.5454
$ws.Range("J132").Value = 3115.0833
$ws.Range("L132").Value = 9345.249899999999
$ws.Range("N132").Value = -14405.2499

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2272.625
$ws.Range("I68").Value = 2572.25
$ws.Range("J68").Value = 1973
$ws.Range("K68").Value = 2572.25
$ws.Range("L68").Value = 1973
$ws.Range("M68").Value = -1823.25
$ws.Range("N68").Value = -3471
$ws.Range("H71").Value = 2272.625
$ws.Range("I71").Value = 2572.25
$ws.Range("J71").Value = 1973
$ws.Range("K71").Value = 12861.25
$ws.Range("L71").Value = 9865
$ws.Range("M71").Value = -9117.25
$ws.Range("N71").Value = -17353
$ws.Range("H124").Value = 54999.668
$ws.Range("J124").Value = 54999.668
$ws.Range("L124").Value = 54999.668
$ws.Range("N124").Value = -64819.668
$ws.Range("H132").Value = 3514.5588
$ws.Range("I132").Value = 2752.1365
$ws.Range("K132").Value = 8256.4095
$ws.Range("M132").Value = -5726.4095

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H126").Value = 2298.3333
$ws.Range("I126").Value = 1748.2142
$ws.Range("K126").Value = 5244.642599999999
$ws.Range("M126").Value = -2774.642599999999
